# Regenerate orders with updated distance/size codes.
# The experiment's distance and size condition labels changed:
#   D51 -> D55, D80 -> D86, D64 -> D69 (distance codes)
#   S30 -> S31                         (size code)
# These codes appear embedded inside many string values across the sheet
# (Condition, Filename_Left, Filename_Right, Distance, Size columns), so we
# do a whole-sheet exact text substitution for each code, oldest-safe order
# doesn't matter since none of the replacement targets collide with any of
# the source codes or each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

[void]$rng.Replace("D51", "D55")
[void]$rng.Replace("D80", "D86")
[void]$rng.Replace("D64", "D69")
[void]$rng.Replace("S30", "S31")
